$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12, shifting existing rows 12-19 down to 13-20
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with the new data record
$ws.Cells.Item(12, 1).Value = 8
$ws.Cells.Item(12, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(12, 3).Value = "Coquimbo"
$ws.Cells.Item(12, 4).Value = 44781
$ws.Cells.Item(12, 4).NumberFormat = $ws.Cells.Item(11, 4).NumberFormat
$ws.Cells.Item(12, 5).Value = 4
$ws.Cells.Item(12, 6).Value = 100112026
$ws.Cells.Item(12, 7).Value = "Haba"
$ws.Cells.Item(12, 8).Value = "Sin especificar"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 400
$ws.Cells.Item(12, 11).Value = 10000
$ws.Cells.Item(12, 12).Value = 11000
$ws.Cells.Item(12, 13).Value = 10500
$ws.Cells.Item(12, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(12, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(12, 16).Value = 420
$ws.Cells.Item(12, 17).Value = 25
$ws.Cells.Item(12, 18).Value = "Hortaliza"
